$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J needs the same header style (bold, centered, bordered) as the
# existing header cells, so copy formatting from I1 before writing into it.
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 10))

# Header row (row 1), including new column J "pete-the-baker"
$headers = @("username", "counting-sheep-dot-dot-dot", "fake-binary", "counting-in-the-amazon", "deodorant-evaporator", "tests-results", "count-the-smiley-faces", "keep-hydrated-1", "see-you-next-happy-year", "pete-the-baker")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Student rows (row 2..27): username, 8 status columns (B..I), new boolean column J (pete-the-baker = No/false for everyone)
$students = @(
    @("a_romegar", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("Afidalgo-fmm", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("AnaWalsh", "Good", "No", "No", "Good", "No", "No", "Good", "Good"),
    @("Andrestart", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("baccandres", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("beatrizsp", "Good", "Good", "No", "Good", "No", "Good", "Good", "No"),
    @("bvispo", "Good", "Good", "Good", "Good", "No", "No", "Good", "No"),
    @("Carlosleono", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("Danihelguera", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("EduardoRivera_98", "Good", "Good", "Good", "Good", "Good", "No", "Good", "Good"),
    @("elliotesp", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "No"),
    @("sinatxe", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("FerZZ", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("Pnat1", "No", "Good", "Good", "Good", "Good", "No", "Good", "Good"),
    @("Juanpda", "No", "Good", "Good", "Good", "Good", "No", "Good", "Good"),
    @("juanchovpa", "Good", "Good", "Good", "No", "No", "Good", "Good", "No"),
    @("LuisSerranoCerame", "Good", "Good", "Good", "No", "Good", "Good", "Good", "Good"),
    @("LydiaAR", "Good", "Good", "No", "Good", "No", "No", "Good", "No"),
    @("mariaperezdeayalas", "Good", "Good", "Good", "Good", "No", "Good", "Good", "No"),
    @("mariadelas", "Good", "Good", "Good", "No", "Good", "Good", "Good", "No"),
    @("NicolasPce", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("pedromartinezleis", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("Roblak96", "No", "Good", "No", "Good", "No", "Good", "Good", "No"),
    @("silviarico", "No", "Good", "Good", "No", "Good", "Good", "Good", "No"),
    @("TheIronhidex", "Good", "Good", "Good", "Good", "Good", "Good", "Good", "Good"),
    @("daniela-arias", "Good", "Good", "No", "Good", "No", "Good", "Good", "Good")
)

for ($i = 0; $i -lt $students.Length; $i++) {
    $row = $students[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $ws.Cells.Item($r, 10).Value = $false
}
